# Update the date in the title paragraph.
$d = $word.ActiveDocument
$d.Content.Find.Execute("2025-12-01 Monday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2025-12-02 Tuesday", 2)

# Update the division problems in the table. The table has 20 rows x 5
# columns, but only every 4th row (1, 5, 9, 13, 17 in 1-based indexing)
# actually holds text; the others are blank spacer rows. Addressing cells
# by row/column (rather than a global find/replace) is required because
# some of the old values ("29÷5=5, 4") repeat but map to different new
# values depending on position.
$t = $d.Tables.Item(1)

$updates = @{
    1  = @("73÷8=9, 1", "97÷6=16, 1", "31÷8=3, 7", "88÷9=9, 7", "68÷5=13, 3")
    5  = @("78÷9=8, 6", "10÷9=1, 1", "46÷9=5, 1", "51÷9=5, 6", "66÷6=11, 0")
    9  = @("94÷2=47, 0", "44÷6=7, 2", "27÷7=3, 6", "50÷9=5, 5", "70÷3=23, 1")
    13 = @("27÷8=3, 3", "48÷9=5, 3", "39÷5=7, 4", "38÷4=9, 2", "53÷8=6, 5")
    17 = @("84÷5=16, 4", "59÷7=8, 3", "46÷9=5, 1", "90÷7=12, 6", "27÷4=6, 3")
}

foreach ($rowIndex in $updates.Keys) {
    $row = $t.Rows.Item($rowIndex)
    $values = $updates[$rowIndex]
    for ($col = 1; $col -le 5; $col++) {
        $row.Cells.Item($col).Range.Text = $values[$col - 1]
    }
}
